$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45185
}
